$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1. Drop the whole draft intro block (from "Nelle marche ci sono X
#    comuni" through the empty paragraph that follows "mappa e tabella
#    con primi 15"). The paragraph that begins with "Presenza di servizi
#    domiciliari " becomes the new first paragraph of the document.
# -----------------------------------------------------------------------
$anchor = $d.Content.Find
$anchor.ClearFormatting()
$anchor.Text = "Presenza di servizi domiciliari"
$found = $anchor.Execute()

if ($found) {
    $keepStart = $anchor.Parent.Start
    $introRange = $d.Range(0, $keepStart)
    $introRange.Delete()
}

# -----------------------------------------------------------------------
# 2. The run that used to open the new page ("EXP Spiegare come abbiamo
#    calcolato ...") carried a stale <w:lastRenderedPageBreak/> marker
#    left over from the previous layout/pagination. Re-issuing its text
#    through Find/Replace rewrites the run and drops the stale marker.
# -----------------------------------------------------------------------
$pageBreakRun = $d.Content.Find
$pageBreakRun.ClearFormatting()
$pageBreakRun.Execute("EXP Spiegare come abbiamo calcolato", $true, $false, $false, $false, $false, $true, 1, $false, "EXP Spiegare come abbiamo calcolato", 2)
